# Apply the rent/mortgage fixture update:
#  - Maintenance Expenses drops from ₹ 3,000 to ₹ 2,500 (-500/month)
#  - Total Monthly Payment drops from ₹ 69,767 to ₹ 69,267 (-500/month), everywhere it is quoted
#  - Every per-period "Taxes, Home Insurance & Maintenance" line item of ₹ 4,250 drops to ₹ 3,750
#  - Yearly/summary rollups of those figures shrink by the same per-period amount

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unique, one-off cells first.
$ws.Range("B5").Value = "₹ 2,500"
$ws.Range("B11").Value = "₹ 4,50,000"
$ws.Range("D13").Value = "₹ 26,250"
$ws.Range("E13").Value = "₹ 4,84,868"
$ws.Range("D148").Value = "₹ 18,750"
$ws.Range("E148").Value = "₹ 3,46,335"

# Global replacements that recur across the monthly schedule (both the
# standalone cells and the concatenated per-month strings in column A).
[void]$ws.Cells.Replace("₹ 4,250", "₹ 3,750")
[void]$ws.Cells.Replace("₹ 69,767", "₹ 69,267")
[void]$ws.Cells.Replace("₹ 51,000", "₹ 45,000")
[void]$ws.Cells.Replace("₹ 8,37,203", "₹ 8,31,203")
